$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1304347826086956
$ws.Range("C2").Value = 0.6086956521739131
$ws.Range("P2").Value = 0.08695652173913043
$ws.Range("S2").Value = 0.1739130434782609
$ws.Range("J3").Value = 0.1428571428571428
$ws.Range("P3").Value = 0.7142857142857143
$ws.Range("S3").Value = 0.1428571428571428
$ws.Range("J4").Value = 0.5
$ws.Range("P4").Value = 0.5
$ws.Range("S5").Value = 1
$ws.Range("B6").Value = 0.1333333333333333
$ws.Range("F6").Value = 0.06666666666666667
$ws.Range("J6").Value = 0.3333333333333333
$ws.Range("O6").Value = 0.06666666666666667
$ws.Range("Q6").Value = 0.1333333333333333
$ws.Range("R6").Value = 0.06666666666666667
$ws.Range("S6").Value = 0.2
$ws.Range("J7").Value = 0.1111111111111111
$ws.Range("Q7").Value = 0.5555555555555556
$ws.Range("R7").Value = 0.1111111111111111
$ws.Range("S7").Value = 0.2222222222222222
$ws.Range("D8").Value = 0.03846153846153846
$ws.Range("E8").Value = 0.03846153846153846
$ws.Range("F8").Value = 0.1153846153846154
$ws.Range("J8").Value = 0.1923076923076923
$ws.Range("O8").Value = 0.03846153846153846
$ws.Range("Q8").Value = 0.1538461538461539
$ws.Range("R8").Value = 0.03846153846153846
$ws.Range("S8").Value = 0.3846153846153846
$ws.Range("B9").Value = 0.1538461538461539
$ws.Range("J9").Value = 0.1538461538461539
$ws.Range("Q9").Value = 0.3076923076923077
$ws.Range("R9").Value = 0.2307692307692308
$ws.Range("S9").Value = 0.1538461538461539
$ws.Range("B10").Value = 0.1322314049586777
$ws.Range("D10").Value = 0.008264462809917356
$ws.Range("F10").Value = 0.05785123966942149
$ws.Range("J10").Value = 0.140495867768595
$ws.Range("O10").Value = 0.01652892561983471
$ws.Range("Q10").Value = 0.1900826446280992
$ws.Range("R10").Value = 0.1322314049586777
$ws.Range("S10").Value = 0.3223140495867768
$ws.Range("G11").Value = 0.1176470588235294
$ws.Range("J11").Value = 0.1176470588235294
$ws.Range("K11").Value = 0.2352941176470588
$ws.Range("L11").Value = 0.5294117647058824
$ws.Range("G12").Value = 0.5
$ws.Range("J12").Value = 0.4
$ws.Range("L12").Value = 0.1
$ws.Range("G13").Value = 1
$ws.Range("H15").Value = 0.1428571428571428
$ws.Range("I15").Value = 0.07142857142857142
$ws.Range("J15").Value = 0.2857142857142857
$ws.Range("M15").Value = 0.07142857142857142
$ws.Range("O15").Value = 0.07142857142857142
$ws.Range("S15").Value = 0.3571428571428572
$ws.Range("H16").Value = 0.08333333333333333
$ws.Range("J16").Value = 0.75
$ws.Range("S16").Value = 0.1666666666666667
$ws.Range("H17").Value = 0.02631578947368421
$ws.Range("I17").Value = 0.07894736842105263
$ws.Range("J17").Value = 0.7105263157894737
$ws.Range("K17").Value = 0.05263157894736842
$ws.Range("O17").Value = 0.05263157894736842
$ws.Range("S17").Value = 0.07894736842105263
$ws.Range("H18").Value = 0.1363636363636364
$ws.Range("I18").Value = 0.04545454545454546
$ws.Range("J18").Value = 0.5454545454545454
$ws.Range("K18").Value = 0.04545454545454546
$ws.Range("O18").Value = 0.09090909090909091
$ws.Range("S18").Value = 0.1363636363636364
$ws.Range("F19").Value = 0.01176470588235294
$ws.Range("H19").Value = 0.2235294117647059
$ws.Range("I19").Value = 0.09411764705882353
$ws.Range("J19").Value = 0.3764705882352941
$ws.Range("K19").Value = 0.1176470588235294
$ws.Range("M19").Value = 0.01176470588235294
$ws.Range("O19").Value = 0.04705882352941176
$ws.Range("S19").Value = 0.1176470588235294
